# Append the new daily log row (row 92) to the bottom of the data table
# on the active sheet, mirroring the existing rows' layout:
#   A: date (stored as plain text, e.g. "2025/10/11")
#   B: day-of-week (text, e.g. Saturday in Japanese)
#   C: time (number)
#   D: ranking (number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 92

# Column A holds a date-like string ("2025/10/11") that must be stored as
# literal TEXT, not auto-converted to a serial date number. Temporarily
# force a text number format before assigning the value, then restore the
# default "Normal" style so the cell ends up with no special formatting
# applied (matching the rest of the data rows).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025/10/11"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = "土"
$ws.Cells.Item($row, 3).Value = 7
$ws.Cells.Item($row, 4).Value = 39
